{"js": "// Agregaci\u00f3n de campo de STATUS_LICENCIA\n// Update the \"reanudaci\u00f3n de labores\" form:\n//  - current post category code/description (table 2: \"CATEGOR\u00cdA ACTUAL\")\n//  - dependency name and project key (table 2)\n//  - provisional occupant's name (table 3: \"PERSONA QUE OCUPA DE MANERA PROVISIONAL LA PLAZA\")\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Table with \"CATEGOR\u00cdA ACTUAL\" / \"DEPENDENCIA\" / \"CLAVE DEL PROYECTO\" rows\nconst datosLaborales = tables.items[1];\n// Table with the provisional occupant's name row\nconst personaProvisional = tables.items[2];\n\n// Helper: replace the full text of a (single-paragraph) table cell while\n// keeping the existing run/paragraph formatting intact.\nasync function setCellText(table, rowIndex, colIndex, newText) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const paragraph = paragraphs.items[0];\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) CATEGOR\u00cdA ACTUAL -> CLAVE: 2S0101A -> 2A0508A\nawait setCellText(datosLaborales, 4, 1, \"2A0508A\");\n\n// 2) CATEGOR\u00cdA ACTUAL -> DESCRIPCI\u00d3N: AUXILIAR 1A -> OFICIAL ADMINISTRATIVO 5A\nawait setCellText(datosLaborales, 4, 3, \"OFICIAL ADMINISTRATIVO 5A\");\n\n// 3) DEPENDENCIA: CONTROL DE REC. HUMANOS Y SUELDOS APLICADOS -> CONTABILIDAD GUBERNAMENTAL\nawait setCellText(datosLaborales, 7, 1, \"CONTABILIDAD GUBERNAMENTAL\");\n\n// 4) CLAVE DEL PROYECTO: 1140020000000000220 -> 1140031490300000120\nawait setCellText(datosLaborales, 8, 1, \"1140031490300000120\");\n\n// 5) Apellido paterno (provisional occupant): \"\" -> SALVADOR\nawait setCellText(personaProvisional, 0, 0, \"SALVADOR\");\n\n// 6) Apellido materno (provisional occupant): VACANTE -> JIMENEZ\nawait setCellText(personaProvisional, 0, 1, \"JIMENEZ\");\n\n// 7) Nombre(s) (provisional occupant): \"\" -> ISIDRO NOE\nawait setCellText(personaProvisional, 0, 2, \"ISIDRO NOE\");\n", "ps1": "# Agregaci\u00f3n de campo de STATUS_LICENCIA\n# Update the \"reanudaci\u00f3n de labores\" form:\n#  - current post category code/description (table 2: \"CATEGOR\u00cdA ACTUAL\")\n#  - dependency name and project key (table 2)\n#  - provisional occupant's name (table 3: \"PERSONA QUE OCUPA DE MANERA PROVISIONAL LA PLAZA\")\n\n$d = $word.ActiveDocument\n\n# Helper: set the text of table($tableIndex) row($rowIndex) cell($cellIndex)\n# (all 1-based, matching Word's COM numbering). The table/row/cell chain is\n# re-resolved from $d fresh on every call so each write lands on the live\n# document instead of a stale cached reference.\nfunction Set-CellText($tableIndex, $rowIndex, $cellIndex, $newText) {\n    $table = $d.Tables.Item($tableIndex)\n    $row = $table.Rows.Item($rowIndex)\n    $cell = $row.Cells.Item($cellIndex)\n    $cell.Range.Text = $newText\n}\n\n# Table 2 = \"2.- DATOS LABORALES\"\n# Row 5 = \"CATEGOR\u00cdA ACTUAL\" | CLAVE | (blank) | DESCRIPCI\u00d3N\n# 1) CLAVE: 2S0101A -> 2A0508A\nSet-CellText 2 5 2 \"2A0508A\"\n# 2) DESCRIPCI\u00d3N: AUXILIAR 1A -> OFICIAL ADMINISTRATIVO 5A\nSet-CellText 2 5 4 \"OFICIAL ADMINISTRATIVO 5A\"\n\n# Row 8 = DEPENDENCIA\n# 3) CONTROL DE REC. HUMANOS Y SUELDOS APLICADOS -> CONTABILIDAD GUBERNAMENTAL\nSet-CellText 2 8 2 \"CONTABILIDAD GUBERNAMENTAL\"\n\n# Row 9 = CLAVE DEL PROYECTO\n# 4) 1140020000000000220 -> 1140031490300000120\nSet-CellText 2 9 2 \"1140031490300000120\"\n\n# Table 3 = \"3.- PERSONA QUE OCUPA DE MANERA PROVISIONAL LA PLAZA\"\n# Row 1 = APELLIDO PATERNO | APELLIDO MATERNO | NOMBRE(S) (values row, above the labels row)\n# 5) Apellido paterno: (blank) -> SALVADOR\nSet-CellText 3 1 1 \"SALVADOR\"\n# 6) Apellido materno: VACANTE -> JIMENEZ\nSet-CellText 3 1 2 \"JIMENEZ\"\n# 7) Nombre(s): (blank) -> ISIDRO NOE\nSet-CellText 3 1 3 \"ISIDRO NOE\"\n"}
